$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 3013426
$ws.Range("I129").Value = 35715510
$ws.Range("J129").Value = 1391.6578
$ws.Range("K129").Value = 107146530
$ws.Range("L129").Value = 4174.9734
$ws.Range("M129").Value = -107141530
$ws.Range("N129").Value = -14174.9734

$ws.Range("H137").Value = 3849803.5
$ws.Range("I137").Value = 8338441
$ws.Range("J137").Value = 2399.8572
$ws.Range("K137").Value = 25015323
$ws.Range("L137").Value = 7199.5716
$ws.Range("M137").Value = -25012773
$ws.Range("N137").Value = -12299.5716

$ws.Range("H138").Value = 3269.9119
$ws.Range("I138").Value = 1517.8667
$ws.Range("J138").Value = 6697.826
$ws.Range("K138").Value = 4553.6001
$ws.Range("L138").Value = 20093.478
$ws.Range("M138").Value = 586.3999000000003
$ws.Range("N138").Value = -30373.478

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4395.87
$ws.Range("I32").Value = 2971.4478
$ws.Range("K32").Value = 2971.4478
$ws.Range("M32").Value = -2684.4478

$ws.Range("H52").Value = 39780
$ws.Range("J52").Value = 39780
$ws.Range("L52").Value = 39780
$ws.Range("N52").Value = -40416

$ws.Range("H97").Value = 591
$ws.Range("I97").Value = 591
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 591
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -95
$ws.Range("N97").ClearContents()

$ws.Range("H132").Value = 12197620
$ws.Range("I132").Value = 15627268
$ws.Range("J132").Value = 3316.5557
$ws.Range("K132").Value = 46881804
$ws.Range("L132").Value = 9949.6671
$ws.Range("M132").Value = -46879274
$ws.Range("N132").Value = -15009.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 30690
$ws.Range("J51").Value = 30690
$ws.Range("L51").Value = 30690
$ws.Range("N51").Value = -31672

$ws.Range("H94").Value = 590.2593000000001
$ws.Range("I94").Value = 471.85
$ws.Range("J94").Value = 928.5714
$ws.Range("K94").Value = 471.85
$ws.Range("L94").Value = 928.5714
$ws.Range("M94").Value = -20.85000000000002
$ws.Range("N94").Value = -1830.5714

$ws.Range("H105").Value = 2062.4707
$ws.Range("I105").Value = 1558.7778
$ws.Range("J105").Value = 2629.125
$ws.Range("K105").Value = 1558.7778
$ws.Range("L105").Value = 2629.125
$ws.Range("M105").Value = 188.2221999999999
$ws.Range("N105").Value = -6123.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1597
$ws.Range("I22").Value = 278.2
$ws.Range("J22").Value = 2915.8
$ws.Range("K22").Value = 278.2
$ws.Range("L22").Value = 2915.8
$ws.Range("M22").Value = 71.80000000000001
$ws.Range("N22").Value = -3615.8

$ws.Range("H25").Value = 29185.2
$ws.Range("I25").Value = 1966.6666
$ws.Range("J25").Value = 70013
$ws.Range("K25").Value = 1966.6666
$ws.Range("L25").Value = 70013
$ws.Range("N25").Value = -70361
$ws.Range("M25").Value = -1792.6666

$ws.Range("H32").Value = 30604.4
$ws.Range("I32").Value = 3000
$ws.Range("J32").Value = 37505.5
$ws.Range("K32").Value = 3000
$ws.Range("L32").Value = 37505.5
$ws.Range("M32").Value = -2684
$ws.Range("N32").Value = -38137.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 2966.6667
$ws.Range("J19").Value = 2966.6667
$ws.Range("L19").Value = 8900.0001
$ws.Range("N19").Value = -9248.0001

$ws.Range("H22").Value = 6667
$ws.Range("I22").Value = 660
$ws.Range("J22").Value = 10957.714
$ws.Range("K22").Value = 1980
$ws.Range("L22").Value = 32873.142
$ws.Range("M22").Value = -1811
$ws.Range("N22").Value = -33211.142

$ws.Range("H27").Value = 6667
$ws.Range("I27").Value = 660
$ws.Range("J27").Value = 10957.714
$ws.Range("K27").Value = 1980
$ws.Range("L27").Value = 32873.142
$ws.Range("M27").Value = -1878
$ws.Range("N27").Value = -33077.142

$ws.Range("H68").Value = 2110.5635
$ws.Range("I68").Value = 698.4
$ws.Range("K68").Value = 2095.2
$ws.Range("M68").Value = -1284.2

$ws.Range("H71").Value = 2110.5635
$ws.Range("I71").Value = 698.4
$ws.Range("K71").Value = 6285.599999999999
$ws.Range("M71").Value = -2229.599999999999

$ws.Range("H80").Value = 5980.6
$ws.Range("J80").Value = 5980.6
$ws.Range("L80").Value = 17941.8
$ws.Range("N80").Value = -19813.8

$ws.Range("H83").Value = 5980.6
$ws.Range("J83").Value = 5980.6
$ws.Range("L83").Value = 53825.4
$ws.Range("N83").Value = -63185.4

$ws.Range("H131").Value = 1787.0857
$ws.Range("I131").Value = 2885.4546
$ws.Range("J131").Value = 1283.6666
$ws.Range("K131").Value = 8656.363799999999
$ws.Range("L131").Value = 3850.9998
$ws.Range("M131").Value = -3616.363799999999
$ws.Range("N131").Value = -13930.9998

$ws.Range("H137").Value = 3070.76
$ws.Range("I137").Value = 3026.0557
$ws.Range("J137").Value = 3185.7144
$ws.Range("K137").Value = 9078.167099999999
$ws.Range("L137").Value = 9557.143199999999
$ws.Range("M137").Value = -3978.167099999999
$ws.Range("N137").Value = -19757.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3140.375
$ws.Range("I126").Value = 1657.9
$ws.Range("J126").Value = 4199.2856
$ws.Range("K126").Value = 4973.700000000001
$ws.Range("L126").Value = 12597.8568
$ws.Range("M126").Value = -2503.700000000001
$ws.Range("N126").Value = -17537.8568

$ws.Range("H132").Value = 3548.4211
$ws.Range("I132").Value = 2590.077
$ws.Range("K132").Value = 7770.231000000001
$ws.Range("M132").Value = -5240.231000000001

$ws.Range("H138").Value = 47313
$ws.Range("J138").Value = 47313
$ws.Range("L138").Value = 47313
$ws.Range("N138").Value = -57593

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 125002400
$ws.Range("I22").Value = 333334340
$ws.Range("J22").Value = 3239.8
$ws.Range("K22").Value = 333334340
$ws.Range("L22").Value = 3239.8
$ws.Range("M22").Value = -333334045
$ws.Range("N22").Value = -3829.8

$ws.Range("H27").Value = 125002400
$ws.Range("I27").Value = 333334340
$ws.Range("J27").Value = 3239.8
$ws.Range("K27").Value = 333334340
$ws.Range("L27").Value = 3239.8
$ws.Range("M27").Value = -333334233
$ws.Range("N27").Value = -3453.8

$ws.Range("H55").Value = 7000.3335
$ws.Range("I55").Value = 5666.6665
$ws.Range("J55").Value = 8334
$ws.Range("K55").Value = 5666.6665
$ws.Range("L55").Value = 8334
$ws.Range("M55").Value = -5493.6665
$ws.Range("N55").Value = -8680

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3705778.2
$ws.Range("I126").Value = 1966.9333
$ws.Range("J126").Value = 8335542.5
$ws.Range("K126").Value = 5900.7999
$ws.Range("L126").Value = 25006627.5
$ws.Range("M126").Value = -3430.7999
$ws.Range("N126").Value = -25011567.5
